# Update EUR->ARS rate: 2025-09-17T15:20:05Z
# Appends a new row (23) to the rate-history sheet with the latest
# quote: date, time, and the EUR->ARS conversion string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as plain text (e.g. "2025-09-04") in every
# existing row. Excel's COM layer auto-parses an ISO "yyyy-mm-dd"
# string into a real date serial on assignment, so force the cell to
# Text format first to keep it a literal string, matching the rest of
# the column.
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "2025-09-17"

# Column B (time-of-day) and C (rate description) are stored as plain
# text already and are not subject to that auto-conversion.
$ws.Range("B23").Value = "15:20:05"
$ws.Range("C23").Value = "1.00 EUR = 1,749.6249"
